$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$ws.Range("A14:B17").Select()
$ws.Range("A14").Activate()
